$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names / links) - swap + refresh values ---
$textUpdates = @{
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'B43' = 'KickToken'
    'C43' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Numeric-looking text cells (price / volume %) - force text format so Excel
#     keeps the literal representation instead of auto-converting to a number ---
$numericTextUpdates = @{
    'D2' = '255.04'
    'E2' = '3.69%'
    'D3' = '27.96'
    'E3' = '-4.39%'
    'D4' = '5.356'
    'D5' = '0.05824'
    'E5' = '0.61%'
    'D6' = '6.711'
    'E7' = '0.83%'
    'D8' = '0.9156'
    'E8' = '5.93%'
    'D9' = '0.1420'
    'E9' = '3.38%'
    'D10' = '0.07175'
    'E10' = '1.44%'
    'D11' = '0.03204'
    'E11' = '-1.70%'
    'D12' = '0.09235'
    'E12' = '-1.35%'
    'D13' = '0.001542'
    'E13' = '0.90%'
    'D14' = '0.0006059'
    'E14' = '-94.08%'
    'D15' = '0.005914'
    'E15' = '-2.84%'
    'D16' = '3.500'
    'E16' = '0.40%'
    'E17' = '1.61%'
    'D18' = '2.253'
    'E18' = '3.88%'
    'E19' = '-1.00%'
    'D20' = '0.03445'
    'E20' = '2.98%'
    'E21' = '1.75%'
    'D22' = '3.524'
    'E22' = '11.12%'
    'D23' = '0.04151'
    'E23' = '0.21%'
    'E24' = '-1.38%'
    'D25' = '0.005114'
    'E25' = '23.50%'
    'D26' = '0.001223'
    'E26' = '-0.22%'
    'E27' = '-0.74%'
    'E28' = '34.23%'
    'D40' = '0.03844'
    'E40' = '2.95%'
    'E41' = '2.81%'
    'D42' = '0.002199'
    'E42' = '0.09%'
    'D43' = '0.002949'
    'E43' = '-49.18%'
    'D44' = '0.01001'
    'E44' = '9.34%'
    'D45' = '0.00005289'
    'E45' = '0.21%'
    'E46' = '0.09%'
    'D47' = '0.09999'
    'E47' = '72.59%'
    'D48' = '0.002210'
    'E48' = '1.72%'
    'D49' = '0.00002099'
    'E49' = '0.09%'
    'E50' = '0.09%'
}
foreach ($addr in $numericTextUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $numericTextUpdates[$addr]
}
